# Updates the "cryptos" price-tracker sheet with refreshed prices and
# 1h volume-change percentages (GitHub Actions scheduled refresh), and
# fixes two pairs of rows where the scraped coin data had shifted by one
# rank (Monero/EthereumClassic and Stellar/WhiteBITCoin swapped places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes a value into a cell while always keeping it stored as literal
# text (column D holds price strings such as "1.00" or "56.860.58" that
# Excel would otherwise silently reinterpret as numbers, which would
# strip trailing zeros / change the representation). Re-applying the
# "Normal" style afterwards removes the temporary text number-format
# again so cells that did not have an explicit style keep none.
function Set-TextValue($cell, $value) {
    $needsTextFormat = ($value -match '^[+-]?[0-9]*\.?[0-9]+$')
    if ($needsTextFormat) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
    if ($needsTextFormat) {
        $cell.Style = "Normal"
    }
}

$updates = @(
    @{ Row=2; D="56.860.58"; E="  -3.75%  " }
    @{ Row=3; D="2.539.97"; E="  -4.61%  " }
    @{ Row=4; E="  -0.02%  " }
    @{ Row=5; D="513.39"; E="  -2.30%  " }
    @{ Row=6; D="140.01"; E="  -3.08%  " }
    @{ Row=7; E="  -0.01%  " }
    @{ Row=8; D="0.554"; E="  -2.74%  " }
    @{ Row=9; D="6.49"; E="  -7.11%  " }
    @{ Row=10; D="0.0989"; E="  -3.94%  " }
    @{ Row=11; D="0.323"; E="  -3.86%  " }
    @{ Row=12; E="  -0.32%  " }
    @{ Row=13; D="2.984.70"; E="  -4.65%  " }
    @{ Row=14; D="56.854.39"; E="  -3.74%  " }
    @{ Row=15; D="20.03"; E="  -5.02%  " }
    @{ Row=16; E="  -3.13%  " }
    @{ Row=17; D="2.523.58"; E="  -6.16%  " }
    @{ Row=18; D="331.97"; E="  -2.00%  " }
    @{ Row=19; D="4.28"; E="  -2.67%  " }
    @{ Row=20; D="10.07"; E="  -3.24%  " }
    @{ Row=21; D="6.12"; E="  -4.37%  " }
    @{ Row=22; D="1.00"; E="  +0.21%  " }
    @{ Row=23; D="64.07"; E="  -0.46%  " }
    @{ Row=24; E="  -0.79%  " }
    @{ Row=25; D="1.00"; E="  +0.11%  " }
    @{ Row=26; D="0.399"; E="  -4.59%  " }
    @{ Row=27; D="2.653.91"; E="  -4.60%  " }
    @{ Row=28; D="6.88"; E="  -3.03%  " }
    @{ Row=29; D="0.0₃0749"; E="  -6.52%  " }
    @{ Row=30; E="  -0.05%  " }
    @{ Row=31; D="6.25"; E="  -6.58%  " }
    @{ Row=32; E="  -2.71%  " }
    @{ Row=33; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="18.44"; E="  -2.24%  " }
    @{ Row=34; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="148.03"; E="  -1.86%  " }
    @{ Row=35; D="3.97"; E="  -4.47%  " }
    @{ Row=36; E="  -5.47%  " }
    @{ Row=37; E="  -5.88%  " }
    @{ Row=38; D="35.49"; E="  -3.85%  " }
    @{ Row=39; D="0.822"; E="  -5.86%  " }
    @{ Row=40; D="1.42"; E="  -2.58%  " }
    @{ Row=41; D="0.999"; E="  -0.02%  " }
    @{ Row=42; D="3.47"; E="  -3.32%  " }
    @{ Row=43; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.0952"; E="  -1.91%  " }
    @{ Row=44; B="WhiteBITCoin"; C="https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D="10.60"; E="  -0.58%  " }
    @{ Row=45; D="0.576"; E="  -6.74%  " }
    @{ Row=46; D="259.70"; E="  -5.76%  " }
    @{ Row=47; D="0.0518"; E="  -2.41%  " }
    @{ Row=48; D="18.42"; E="  -7.49%  " }
    @{ Row=49; D="1.964.78"; E="  -4.02%  " }
    @{ Row=50; D="0.0220"; E="  -4.08%  " }
    @{ Row=51; D="4.51"; E="  -3.99%  " }
)

foreach ($update in $updates) {
    $row = $update.Row
    foreach ($col in @("B", "C", "D", "E")) {
        if ($update.ContainsKey($col)) {
            Set-TextValue $ws.Range("$col$row") $update[$col]
        }
    }
}
